$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "91.651.02"
$ws.Range("E2").Value = "  +0.92%  "

# Row 3
$ws.Range("D3").Value = "3.119.00"
$ws.Range("E3").Value = "  +1.25%  "

# Row 4
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").Value = "'246.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.04%  "

# Row 6
$ws.Range("D6").Value = "'617.42"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.15%  "

# Row 7
$ws.Range("E7").Value = "  -2.23%  "

# Row 8
$ws.Range("E8").Value = "  +4.89%  "

# Row 9
$ws.Range("E9").Value = "  -0.16%  "

# Row 10
$ws.Range("D10").Value = "3.117.04"
$ws.Range("E10").Value = "  +1.35%  "

# Row 11
$ws.Range("D11").Value = "'0.734"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.82%  "

# Row 12
$ws.Range("E12").Value = "  +1.68%  "

# Row 13
$ws.Range("E13").Value = "  +1.13%  "

# Row 14
$ws.Range("E14").Value = "  +3.15%  "

# Row 15
$ws.Range("E15").Value = "  -0.63%  "

# Row 16
$ws.Range("D16").Value = "91.571.44"
$ws.Range("E16").Value = "  +0.58%  "

# Row 17
$ws.Range("D17").Value = "3.701.49"

# Row 18
$ws.Range("D18").Value = "3.083.95"
$ws.Range("E18").Value = "  -0.42%  "

# Row 19
$ws.Range("E19").Value = "  +0.51%  "

# Row 20
$ws.Range("D20").Value = "'14.86"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.00%  "

# Row 21
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'9.53"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.66%  "

# Row 22
$ws.Range("B22").Value = "Polkadot"
$ws.Range("C22").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D22").Value = "'5.81"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.50%  "

# Row 23
$ws.Range("D23").Value = "'446.50"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.42%  "

# Row 24
$ws.Range("D24").Value = "'0.0000202"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.38%  "

# Row 25
$ws.Range("D25").Value = "'5.85"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.82%  "

# Row 26
$ws.Range("D26").Value = "'87.94"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.29%  "

# Row 27
$ws.Range("D27").Value = "'11.75"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.04%  "

# Row 28
$ws.Range("D28").Value = "3.282.45"
$ws.Range("E28").Value = "  +0.68%  "

# Row 29
$ws.Range("D29").Value = "'0.145"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +31.04%  "

# Row 30
$ws.Range("E30").Value = "  +0.19%  "

# Row 31
$ws.Range("D31").Value = "'0.235"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.10%  "

# Row 32
$ws.Range("D32").Value = "'0.168"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -8.68%  "

# Row 33
$ws.Range("E33").Value = "  +4.54%  "

# Row 34
$ws.Range("D34").Value = "'9.30"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.91%  "

# Row 35
$ws.Range("E35").Value = "  -1.04%  "

# Row 36
$ws.Range("D36").Value = "'7.88"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.50%  "

# Row 37
$ws.Range("D37").Value = "'26.21"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.36%  "

# Row 38
$ws.Range("D38").Value = "'4.16"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.68%  "

# Row 39
$ws.Range("E39").Value = "  +1.38%  "

# Row 40
$ws.Range("D40").Value = "'490.54"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.03%  "

# Row 41
$ws.Range("E41").Value = "  +1.56%  "

# Row 42
$ws.Range("D42").Value = "'0.439"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.76%  "

# Row 43
$ws.Range("D43").Value = "'3.41"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.69%  "

# Row 44
$ws.Range("D44").Value = "'22.19"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.31%  "

# Row 46
$ws.Range("D46").Value = "'157.73"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.42%  "

# Row 47
$ws.Range("D47").Value = "'0.707"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.72%  "

# Row 48
$ws.Range("D48").Value = "'1.91"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.33%  "

# Row 49
$ws.Range("E49").Value = "  +2.01%  "

# Row 50
$ws.Range("D50").Value = "'44.05"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.01%  "

# Row 51
$ws.Range("E51").Value = "  -1.38%  "
